$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("B1").Value = "Text"

# Strip trailing " ?"/"?" from each question text (Text column)
$ws.Range("B2").Value = "What is Supervised learning"
$ws.Range("B3").Value = "What is regression"
$ws.Range("B4").Value = "What is slope"
$ws.Range("B5").Value = "What is intercept"
$ws.Range("B6").Value = "What is semi supervised learning "
$ws.Range("B7").Value = "What is Data split "
$ws.Range("B8").Value = "What do you mean by the training of data"
$ws.Range("B9").Value = "What do you mean by the testing of data"
$ws.Range("B10").Value = "What is the Validation of model"
$ws.Range("B11").Value = "What is Feature Selection"
$ws.Range("B12").Value = "What is Feature scaling"
$ws.Range("B13").Value = "What is Model evalution"
$ws.Range("B14").Value = "What is Curse of Dimensionality"
$ws.Range("B15").Value = "How specific output links to Supervised learning"
$ws.Range("B16").Value = "What is confusion matrix"
$ws.Range("B17").Value = "What is Classification"
$ws.Range("B18").Value = "What is Linear Regression"
$ws.Range("B19").Value = "What is Logistic Regression"
$ws.Range("B20").Value = "What is Naïve bayes algorithm"
$ws.Range("B21").Value = "What is euclidean distance"
$ws.Range("B22").Value = "What is Manhattan distance"
$ws.Range("B23").Value = "What is KNN algorithm"
$ws.Range("B24").Value = "What is Mean Squared Error or MSE"
$ws.Range("B25").Value = "What is Accuracy"
$ws.Range("B26").Value = "What is ROC or AUC curve"
$ws.Range("B27").Value = "What is OLS methods"
$ws.Range("B28").Value = "What is Gradient descent"
$ws.Range("B29").Value = "What is Cost Function"
$ws.Range("B30").Value = "What is F1 score"
$ws.Range("B31").Value = "What is Precision"
$ws.Range("B32").Value = "What is Recall"
$ws.Range("B33").Value = "What is True Positive Rate or TPR"
$ws.Range("B34").Value = "What is True Negative Rate or TNR"
$ws.Range("B35").Value = "What is dependent or Target Variable "
$ws.Range("B36").Value = "What is independent or Predictor Variable"
$ws.Range("B37").Value = "What is train test split"
$ws.Range("B38").Value = "How multicolinearity affect the model performance"
$ws.Range("B39").Value = "What is Pearson's correalation coefficient"
$ws.Range("B40").Value = "What is assumption of Linear Regression"
$ws.Range("B41").Value = "What is Sum Of Squared Error or SSE"
$ws.Range("B42").Value = "What is Regression  error or SSR"
$ws.Range("B43").Value = "What is Toatal error or SST"
$ws.Range("B44").Value = "What is Coeeficient Of determinant"
$ws.Range("B45").Value = "What is Adjusted R-square"
$ws.Range("B46").Value = "What is decision Boundary"
$ws.Range("B47").Value = "What is Log loss"
$ws.Range("B48").Value = "What is predict proba"
$ws.Range("B49").Value = "What is imbalance class problem"
$ws.Range("B50").Value = "What  is Polynomial regression"
$ws.Range("B51").Value = "What is Bias Variance trade off"
$ws.Range("B52").Value = "What is Standadization"
$ws.Range("B53").Value = "What is Normalization"
$ws.Range("B54").Value = "What is Posterior probability"
$ws.Range("B55").Value = "What is Predictiving modelling"

# Reset view: clear frozen/scrolled top-left cell and move selection to B2
$ws.Range("B2").Select()

Write-Output "done"
